# "Monthly Feed and Prices"
# The "Whole" sheet (second sheet / ActiveSheet) gets:
#  - A1: the year 2021
#  - Rows 3 (Feed Mass), 4 (Feed Price), 5 (Misc) filled in with monthly 0s
#    and a couple of real values (F3=100, F4=644)
#  - Three new summary rows: Average Age / Feed per Pig / Feed per Pig per age
#  - Column A widened, selection moved to F6

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Whole")

# Year label in the top-left corner of the month table
$ws.Range("A1").Value = 2021

# Row 3 - Feed Mass: fill B:M with 0, then set the May (F) figure
$ws.Range("B3:M3").Value = 0
$ws.Range("F3").Value = 100

# Row 4 - Feed Price: fill B:M with 0, then set the May (F) figure
$ws.Range("B4:M4").Value = 0
$ws.Range("F4").Value = 644

# Row 5 - Misc: fill B:M with 0
$ws.Range("B5:M5").Value = 0

# New rows for the monthly feed/age summary
$ws.Range("A6").Value = "Average Age"
$ws.Range("C6").Value = ""
$ws.Range("A7").Value = "Feed per Pig"
$ws.Range("A8").Value = "Feed per Pig per age"

# Column A is widened to fit the new longer labels
$ws.Columns.Item(1).ColumnWidth = 20.74

# Move the active selection to F6, matching the new working area
$ws.Range("F6").Select()
